$d = $word.ActiveDocument

# The requested change rewrites word/styles.xml's <w:docDefaults> block,
# stripping it down to only the properties that survive in the target
# (font family, sz/szCs, lang for rPrDefault; a bare line-spacing for
# pPrDefault). That block isn't reachable through the normal Styles
# collection (writes to Styles("Normal") land on the Normal style's own
# w:pPr/w:rPr, not on docDefaults), so we edit the underlying OOXML
# package text directly via Document.WordOpenXML (a flat-OPC dump of the
# whole package) and write it back.

$oldDocDefaults = '<w:docDefaults><w:rPrDefault><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b w:val="0"/><w:i w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/><w:shd w:val="clear" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:lang w:val="en"/></w:rPr></w:rPrDefault><w:pPrDefault><w:pPr><w:keepNext w:val="0"/><w:keepLines w:val="0"/><w:widowControl/><w:pBdr><w:top w:val="nil" w:sz="0" w:space="0"/><w:left w:val="nil" w:sz="0" w:space="0"/><w:bottom w:val="nil" w:sz="0" w:space="0"/><w:right w:val="nil" w:sz="0" w:space="0"/><w:between w:val="nil" w:sz="0" w:space="0"/></w:pBdr><w:shd w:val="clear" w:fill="auto"/><w:spacing w:before="0" w:after="0" w:line="276" w:lineRule="auto"/><w:ind w:left="0" w:right="0" w:firstLine="0"/><w:contextualSpacing w:val="0"/><w:jc w:val="left"/></w:pPr></w:pPrDefault></w:docDefaults>'

$newDocDefaults = '<w:docDefaults><w:rPrDefault><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en"/></w:rPr></w:rPrDefault><w:pPrDefault><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr></w:pPrDefault></w:docDefaults>'

$packageXml = $d.WordOpenXML

if ($packageXml.Contains($oldDocDefaults)) {
    $packageXml = $packageXml.Replace($oldDocDefaults, $newDocDefaults)
    $d.WordOpenXML = $packageXml
    Write-Output "docDefaults updated"
} else {
    Write-Output "docDefaults block not found verbatim; no change made"
}
